# Implement recall block and draw power, add 4 new cards
#
# - Adds a "Done" tracking column (boolean checkbox) to the Skill and
#   Power tables (the Attack table already has one).
# - Marks a handful of existing Attack/Skill cards as Done.
# - Names several previously-blank card rows across Attack/Skill.
# - Renames "For Loop" -> "Offensive Loop" (Attack) and
#   "While Loop" -> "Defensive Loop" (Skill), and moves "Rewinder" to a
#   new row, naming its old slot "Accelerate".
# - Leaves the active sheet on "Skill".

$wb = $excel.ActiveWorkbook

$wsRelic  = $wb.Worksheets.Item("Relic")
$wsAttack = $wb.Worksheets.Item("Attack")
$wsSkill  = $wb.Worksheets.Item("Skill")
$wsPower  = $wb.Worksheets.Item("Power")

# ---------------------------------------------------------------------
# Attack sheet ("Table1") already has its "Done" column (col I) - just
# fill in the new data / names.
# ---------------------------------------------------------------------
$wsAttack.Range("I4").Value = $true

$wsAttack.Range("A17").Value = "Flash Forward"
$wsAttack.Range("A20").Value = "Flashback"
$wsAttack.Range("A21").Value = "Offensive Loop"
$wsAttack.Range("A29").Value = "Astral Banishment"

$wsAttack.Range("E4").Select()

# ---------------------------------------------------------------------
# Skill sheet ("Table13") - add the "Done" column by resizing the table
# and naming the new header, then populate it.
# ---------------------------------------------------------------------
$loSkill = $wsSkill.ListObjects.Item(1)
$loSkill.Resize($wsSkill.Range("A1:I31"))

$wsSkill.Range("H1").Copy()
$wsSkill.Range("I1").PasteSpecial(-4122) # xlPasteFormats
$wsSkill.Range("I1").Value = "Done"

$wsSkill.Range("I2").Value = $true
$wsSkill.Range("I4").Value = $true
$wsSkill.Range("I7").Value = $true

$wsSkill.Range("A3").Value = "Lucid Dream"
$wsSkill.Range("A4").Value = "Mystic Barrier"
$wsSkill.Range("A7").Value = "Delayed Guard"
$wsSkill.Range("A9").Value = "Accelerate"
$wsSkill.Range("A12").Value = "Dark Matter"
$wsSkill.Range("A22").Value = "Defensive Loop"
$wsSkill.Range("A23").Value = "Timestream Shift"
$wsSkill.Range("A24").Value = "Big Bang"
$wsSkill.Range("A25").Value = "Rewinder"
$wsSkill.Range("A31").Value = "Wheel of Time"

$wsSkill.Columns.Item(9).ColumnWidth = 16

# ---------------------------------------------------------------------
# Power sheet ("Table134") - add the (still empty) "Done" column too.
# ---------------------------------------------------------------------
$loPower = $wsPower.ListObjects.Item(1)
$loPower.Resize($wsPower.Range("A1:I13"))

$wsPower.Range("H1").Copy()
$wsPower.Range("I1").PasteSpecial(-4122) # xlPasteFormats
$wsPower.Range("I1").Value = "Done"

$wsPower.Columns.Item(9).ColumnWidth = 16

$wsPower.Range("I3").Select()

# ---------------------------------------------------------------------
# Cosmetic selection tweak on Relic, carried over from the session.
# ---------------------------------------------------------------------
$wsRelic.Range("C2").Select()

# ---------------------------------------------------------------------
# Leave the Skill sheet active/selected, matching the saved workbook
# view state (activeTab moves from Attack to Skill).
# ---------------------------------------------------------------------
$wsSkill.Activate()
$wsSkill.Range("I5").Select()
